$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 13 (old rows 13-23 shift down to 15-25)
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# Row 10
$ws.Cells.Item(10, 1).Value = "Objetivos:"
$ws.Cells.Item(10, 2).Value = "Como parte fundamental da formação específica e geral, a disciplina tem por objetivos (a) fornecer os conceitos fundamentais sobre sistemas de automação da produção e suas características, aplicações, capacidades e princípios de programação, que são importantes para a formação em engenharia física; (b) capacitar o aluno, trabalhando individualmente e em grupo, a modelar e resolver problemas de interesse envolvendo técnicas de automação de manufatura, com escolhas adequadas de hipóteses e aplicação de ferramentas correspondentes de solução; e (c) aplicar e estender os conceitos físicos aprendidos previamente."
$ws.Cells.Item(10, 3).Value = "Como parte fundamental da formação específica e geral, a disciplina tem por objetivos (a) fornecer os conceitos fundamentais sobre sistemas de automação da produção e suas características, aplicações, capacidades e princípios de programação, que são importantes para a formação em engenharia física; (b) capacitar o aluno, trabalhando individualmente e em grupo, a modelar e resolver problemas de interesse envolvendo técnicas de automação de manufatura, com escolhas adequadas de hipóteses e aplicação de ferramentas correspondentes de solução; e (c) aplicar e estender os conceitos físicos aprendidos previamente."
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Cells.Item(11, 1).Value = "Objectives:"
$ws.Cells.Item(11, 2).Value = "As an important requirement for the specific and the general engineering studies, the course aims to (a) provide the fundamental concepts concerning automated manufacturing systems, as well as their applications, characteristics, yield, and programming, which are important for the learning process in physics engineering; (b) enable the student, to work as individually as in groups, to model and solve problems of interest regarding automated manufacturing systems, including the adequate choice of hypotheses and the application of suitable solution tools; and (c) apply and extend the previously learned physical concepts."
$ws.Cells.Item(11, 3).Value = "As an important requirement for the specific and the general engineering studies, the course aims to (a) provide the fundamental concepts concerning automated manufacturing systems, as well as their applications, characteristics, yield, and programming, which are important for the learning process in physics engineering; (b) enable the student, to work as individually as in groups, to model and solve problems of interest regarding automated manufacturing systems, including the adequate choice of hypotheses and the application of suitable solution tools; and (c) apply and extend the previously learned physical concepts."
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Cells.Item(12, 1).Value = "Docentes responsáveis:"
$ws.Rows.Item(12).EntireRow.AutoFit()

# Row 13
$ws.Cells.Item(13, 2).Value = "519033 - Carlos Yujiro Shigue"
$ws.Cells.Item(13, 3).Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows.Item(13).EntireRow.AutoFit()

# Row 14
$ws.Cells.Item(14, 2).Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Cells.Item(14, 3).Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Rows.Item(14).EntireRow.AutoFit()

# Row 15
$ws.Cells.Item(15, 1).Value = "Programa resumido:"
$ws.Cells.Item(15, 2).Value = "Introdução aos sistemas de manufatura; Tecnologias aplicadas à automação de processos industriais; Equipamentos e programação de hardware para supervisão e automação de processos"
$ws.Cells.Item(15, 3).Value = "Introdução aos sistemas de manufatura; Tecnologias aplicadas à automação de processos industriais; Equipamentos e programação de hardware para supervisão e automação de processos"
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Cells.Item(16, 1).Value = "Short syllabus:"
$ws.Cells.Item(16, 2).Value = "Introduction to manufacturing systems; Technologies applied to industrial process automation; Supervision and automation equipment and hardware programming."
$ws.Cells.Item(16, 3).Value = "Introduction to manufacturing systems; Technologies applied to industrial process automation; Supervision and automation equipment and hardware programming."
$ws.Rows.Item(16).RowHeight = 60

# Row 17
$ws.Cells.Item(17, 1).Value = "Programa:"
$ws.Cells.Item(17, 2).Value = "1. Introdução aos sistemas de manufatura;    2. Indústria 4.0;    3. Inteligencia artifical (ai), internet das coisas (IoT) e sistemas ciberfísicos;    4. Sistemas de controle industrial; conceito de PLM e integração com os sistemas de gestão;    5. Monitoramento e supervisão de processos de produção. Sistemas de controle da produção, manufatura sustentável;    6. Componentes de hardware para automação de processos: controle numérico, programação CNC, controle discreto utilizando controladores lógico programáveis e sistemas on-chip;    7. Robótica industrial – programação de robôs e robôs colaborativos;    8. Sistemas de transporte de materiais e sistemas de armazenamento;    9. Identificação automática e captura de dados – tecnologias de inspeção."
$ws.Cells.Item(17, 3).Value = "1. Introdução aos sistemas de manufatura;    2. Indústria 4.0;    3. Inteligencia artifical (ai), internet das coisas (IoT) e sistemas ciberfísicos;    4. Sistemas de controle industrial; conceito de PLM e integração com os sistemas de gestão;    5. Monitoramento e supervisão de processos de produção. Sistemas de controle da produção, manufatura sustentável;    6. Componentes de hardware para automação de processos: controle numérico, programação CNC, controle discreto utilizando controladores lógico programáveis e sistemas on-chip;    7. Robótica industrial – programação de robôs e robôs colaborativos;    8. Sistemas de transporte de materiais e sistemas de armazenamento;    9. Identificação automática e captura de dados – tecnologias de inspeção."
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Cells.Item(18, 1).Value = "Syllabus:"
$ws.Cells.Item(18, 2).Value = "1. Introduction to manufacturing systems;2. 4.0 industry;3. Artificial intelligence (AI), internet of things (IoT), and cyber-physical systems;4. Industrial control systems; PLM concept and integration with management systems;5. Production process monitoring and supervision. Production control systems, sustainable manufacturing;6. Hardware components for process automation: numerical control, CNC programming, discrete control using programmable logic controllers and on-chip systems;7. Industrial robotics – programming of robots and collaborative robots;8. Material transportation systems and storage systems;9. Automatic identification and data capture – inspection technologies."
$ws.Cells.Item(18, 3).Value = "1. Introduction to manufacturing systems;2. 4.0 industry;3. Artificial intelligence (AI), internet of things (IoT), and cyber-physical systems;4. Industrial control systems; PLM concept and integration with management systems;5. Production process monitoring and supervision. Production control systems, sustainable manufacturing;6. Hardware components for process automation: numerical control, CNC programming, discrete control using programmable logic controllers and on-chip systems;7. Industrial robotics – programming of robots and collaborative robots;8. Material transportation systems and storage systems;9. Automatic identification and data capture – inspection technologies."
$ws.Rows.Item(18).RowHeight = 120

# Row 19
$ws.Cells.Item(19, 1).Value = "Avaliação:"
$ws.Rows.Item(19).EntireRow.AutoFit()

# Row 20
$ws.Cells.Item(20, 1).Value = "Método:"
$ws.Cells.Item(20, 2).Value = "Aulas expositivas e práticas de laboratório com interações em grupo para a solução de problemas."
$ws.Cells.Item(20, 3).Value = "Aulas expositivas e práticas de laboratório com interações em grupo para a solução de problemas."
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Cells.Item(21, 1).Value = "Critério:"
$ws.Cells.Item(21, 2).Value = "Média aritmética (M) de provas individuais (P1 e P2) e trabalhos em grupo ao longo do semestre (T), tal que M = 0,3*P1+0,3*P2+0.4*T"
$ws.Cells.Item(21, 3).Value = "Média aritmética (M) de provas individuais (P1 e P2) e trabalhos em grupo ao longo do semestre (T), tal que M = 0,3*P1+0,3*P2+0.4*T"
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Cells.Item(22, 1).Value = "Norma de recuperação:"
$ws.Cells.Item(22, 2).Value = "Média aritmética (M) de uma prova individual (P1) e trabalhos em grupo ao longo do semestre (T), tal que M = 0,5*P1+0.5*T"
$ws.Cells.Item(22, 3).Value = "Média aritmética (M) de uma prova individual (P1) e trabalhos em grupo ao longo do semestre (T), tal que M = 0,5*P1+0.5*T"
$ws.Rows.Item(22).RowHeight = 60

# Row 23
$ws.Cells.Item(23, 1).Value = "Bibliografia:"
$ws.Cells.Item(23, 2).Value = "GROOVER, M.P. Automação Industrial e Sistemas de Manufatura, 561p., 3a Edição - São Paulo, Pearson Prentice Hall, 2011.RICHARD, L.S., ERNEST, L.H., Handbook of Industrial Automation, Marcel Dekker,Inc. NewYork, 2000.ADALBERTO FILHO ET. AL, Automação & Sociedade: Quarta revolução Industrial, um olhar para o Brasil, 1a Edição, Brasport Livros e Multimídia Limitada."
$ws.Cells.Item(23, 3).Value = "GROOVER, M.P. Automação Industrial e Sistemas de Manufatura, 561p., 3a Edição - São Paulo, Pearson Prentice Hall, 2011.RICHARD, L.S., ERNEST, L.H., Handbook of Industrial Automation, Marcel Dekker,Inc. NewYork, 2000.ADALBERTO FILHO ET. AL, Automação & Sociedade: Quarta revolução Industrial, um olhar para o Brasil, 1a Edição, Brasport Livros e Multimídia Limitada."
$ws.Rows.Item(23).RowHeight = 120

# Row 24
$ws.Cells.Item(24, 1).Value = "Requisitos:"
$ws.Rows.Item(24).EntireRow.AutoFit()

# Row 25
$ws.Cells.Item(25, 2).Value = "LOM3264 -  Fundamentos de Controle  (Requisito)`n"
$ws.Cells.Item(25, 3).Value = "LOM3264 -  Fundamentos de Controle  (Requisito)`n"
$ws.Rows.Item(25).RowHeight = 30
